$wb = $excel.ActiveWorkbook

function Get-HyperlinkAddress($ws, $cellAddr) {
    $found = $null
    foreach ($hl in $ws.Hyperlinks) {
        $ref = $hl.Range.Address()
        if ($ref -eq $cellAddr) {
            $found = $hl.Address
        }
    }
    return $found
}

function Set-HandbackRow($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column C: Status -> handed back, now in sync with en-US
    $ws.Range("C2").Value = "Handed back: in sync with en-US"

    # Reuse the source-file hyperlink address (same target as A2) for the
    # new "Latest Target File" column (F2).
    $mdAddress = Get-HyperlinkAddress $ws "`$A`$2"
    $ws.Range("F2").Value = $ws.Range("A2").Value2
    if ($mdAddress) {
        $ws.Hyperlinks.Add($ws.Range("F2"), $mdAddress, [Type]::Missing, [Type]::Missing, $ws.Range("F2").Value2) | Out-Null
    }

    # Reuse the translated-file hyperlink address (same target as D2) for
    # the new "Latest Handback File" column (G2).
    $xlfAddress = Get-HyperlinkAddress $ws "`$D`$2"
    $ws.Range("G2").Value = $ws.Range("D2").Value2
    if ($xlfAddress) {
        $ws.Hyperlinks.Add($ws.Range("G2"), $xlfAddress, [Type]::Missing, [Type]::Missing, $ws.Range("G2").Value2) | Out-Null
    }

    # Column H: Latest Handback DateTime
    $ws.Range("H2").Value = $handbackDateTime
}

Set-HandbackRow "zh-cn" "2016-03-20 14:35:29"
Set-HandbackRow "de-de" "2016-03-20 14:35:34"
